# GuideQuestTable 에 viewInBattle|Bool 컬럼 추가
# Insert a new "viewInBattle|Bool" boolean column right before the existing
# "needCount|Int" column (old column E) on the GuideQuestTable sheet, and
# mark that sheet as the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GuideQuestTable")

# Shift the old column E (and everything after it) one column to the right,
# opening up a blank column E for the new field.
$ws.Columns.Item(5).Insert()

# Header for the new column.
$ws.Range("E1").Value = "viewInBattle|Bool"

# Every existing data row gets TRUE for the new boolean field.
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $true
$ws.Range("E5").Value = $true
$ws.Range("E6").Value = $true
$ws.Range("E7").Value = $true
$ws.Range("E8").Value = $true

# GuideQuestTable becomes the active/selected sheet in the saved workbook.
$ws.Activate()
